# [Kadastro App] Yeni kayit eklendi: 2956
#
# Appends one new record (Kayit No 2956) as row 47 to the two worksheets
# that track this data: "Kayitlar" (the master log, sheet1.xml) and
# "Erdemli" (the per-birim sheet, sheet8.xml). Both sheets mirror the
# same rows, so the same row is written to each.
#
# Kayit No / Parsel Sayisi / Tarih are entered with a leading apostrophe
# so Excel stores them as literal text (matching every other row in the
# column, e.g. "2955", "1", "2025-09-09") instead of re-typing them as a
# number / date.

$wb = $excel.ActiveWorkbook

$newRecord = @{
    KayitNo       = "2956"
    Tarih         = "2025-09-09"
    Birim         = "Erdemli"
    ParselSayisi  = "1"
    Is            = "ÇAP"
    Personeller   = "CEMAL TİMUROĞLU (K.Teknisyeni)"
}

$targetSheets = @("Kayitlar", "Erdemli")

foreach ($sheetName in $targetSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $newRow = $ws.Cells.Item($ws.UsedRange.Rows.Count, 1).Row + 1

    $ws.Cells.Item($newRow, 1).Value = "'" + $newRecord.KayitNo
    $ws.Cells.Item($newRow, 2).Value = "'" + $newRecord.Tarih
    $ws.Cells.Item($newRow, 3).Value = $newRecord.Birim
    $ws.Cells.Item($newRow, 4).Value = "'" + $newRecord.ParselSayisi
    $ws.Cells.Item($newRow, 5).Value = $newRecord.Is
    $ws.Cells.Item($newRow, 6).Value = $newRecord.Personeller
}
